# Weekly fruit/vegetable price update: a new weekly record is inserted as
# row 81, pushing the existing rows 81-128 down to 82-129.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 81 (shifts rows 81.. downward).
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with the new weekly price record.
$ws.Cells.Item(81, 1).Value  = 1
$ws.Cells.Item(81, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(81, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(81, 4).Value  = 44813
$ws.Cells.Item(81, 5).Value  = 15
$ws.Cells.Item(81, 6).Value  = "Fruta"
$ws.Cells.Item(81, 7).Value  = 100106
$ws.Cells.Item(81, 8).Value  = "Oleaginosos"
$ws.Cells.Item(81, 9).Value  = 100106002
$ws.Cells.Item(81, 10).Value = "Palta"
$ws.Cells.Item(81, 11).Value = "Hass"
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 400
$ws.Cells.Item(81, 14).Value = 22000
$ws.Cells.Item(81, 15).Value = 23000
$ws.Cells.Item(81, 16).Value = 22500
$ws.Cells.Item(81, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(81, 18).Value = "Perú"
$ws.Cells.Item(81, 19).Value = 2250
$ws.Cells.Item(81, 20).Value = 10
